# "excel parser fix: empty strings to nulls"
#
# The lessonList sheet had a literal text header row (Subject / Teacher /
# StudentGroup) in row 1 that confused a downstream parser expecting null
# for empty cells. Remove that header row entirely so the real lesson
# data (previously rows 2-21) shifts up to become rows 1-20, then
# reselect A1 so the sheet opens with the first data row highlighted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lessonList")

$ws.Rows("1:1").Delete() | Out-Null
$ws.Range("A1").Select() | Out-Null
